$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove every existing hyperlink on the sheet (Range.Hyperlinks.Delete() acts
# sheet-wide in this host, so one call clears A2/B2/A3/B4's links in one go).
$ws.Range("A1").Hyperlinks.Delete()
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Style = "Normal"

# --- Row 2: shubhamjyani2@gmail.com | password | Login Unsuccessful | Login Unsuccessful ---
$ws.Range("B2").Value = "password"
$ws.Range("C2").Value = "Login Unsuccessful"
$ws.Range("D2").Value = "Login Unsuccessful"

# --- A2's hyperlink (unchanged target, first relationship id) ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:shubhamjyani2@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

# --- Row 3: id@email.com | a3jKkxQB | Login Unsuccessful | Login Unsuccessful ---
# B3 keeps pointing at the old "Kamla@29" mailto, but now displays "a3jKkxQB".
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Kamla@29", "", "", "Kamla@29")
$ws.Range("B3").Value = "a3jKkxQB"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Value = "Login Unsuccessful"
$ws.Range("D3").Value = "Login Unsuccessful"

# --- Row 4: empty | empty | Login Unsuccessful with email and password required error | ... ---
$ws.Range("A4").Value = "empty"
$ws.Range("B4").Value = "empty"
$ws.Range("C4").Value = "Login Unsuccessful with email and password required error"
$ws.Range("D4").Value = "Login Unsuccessful with email and password required error"

# --- A3's hyperlink, text updated last ---
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:shubhamjyani2@gmail.com")
$ws.Range("A3").Value = "id@email.com"
$ws.Range("A3").Style = "Hyperlink"

# Leave selection on A1 (matches the saved file no longer pinning A3).
$ws.Range("A1").Select() | Out-Null
